$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 8 that mirrors row 2 (same scenario), but using a
# dedicated "200" status code (as a shared string) instead of the
# hard-coded numeric 201, and without the ExpectedMessage column.
$ws.Range("A2").Copy($ws.Range("A8"))
$ws.Range("B2").Copy($ws.Range("B8"))
$ws.Range("C2").Copy($ws.Range("C8"))
$ws.Range("D2").Copy($ws.Range("D8"))
$ws.Range("E2").Copy($ws.Range("E8"))
$ws.Range("F2").Copy($ws.Range("F8"))
$ws.Range("G2").Copy($ws.Range("G8"))
$ws.Range("H2").Copy($ws.Range("H8"))
$ws.Range("I2").Copy($ws.Range("I8"))
$ws.Range("J2").Copy($ws.Range("J8"))
$ws.Range("K8").Value = "200"

# Recreate the hyperlink on D8 (mailto link), matching D2's target.
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:orer79521@gmail.com")
# Re-apply D2's (Hyperlink) cell style, since adding the hyperlink above
# resets/duplicates the cell formatting.
$ws.Range("D2").Copy($ws.Range("D8"))

# Update the active selection as reflected in the saved view state.
$ws.Range("E2").Select()
